$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Turn the e-mail address into a mailto: hyperlink, then append a new
#    " | Site: impedrof.github.io/curriculoOn/" segment (also a hyperlink)
#    to the "Dados pessoais" e-mail line.
# ---------------------------------------------------------------------------

# Locate the plain e-mail text and remember where it ends so we can append
# the new runs right after it (still inside the same paragraph).
$emailTail = $d.Content
$emailTail.Find.Execute("pedro.freitas14@fatec.sp.gov.br", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$emailTail.Collapse(0)

# " | "
$pipeStart = $emailTail.Start
$emailTail.InsertAfter(" | ")
$pipeRange = $d.Range($pipeStart, $emailTail.End)
$pipeRange.Font.Size = 11

# "Site:" (bold label)
$siteStart = $emailTail.End
$siteIns = $d.Range($siteStart, $siteStart)
$siteIns.InsertAfter("Site:")
$siteRange = $d.Range($siteStart, $siteIns.End)
$siteRange.Font.Size = 11
$siteRange.Font.Bold = $true

# " "
$spaceStart = $siteIns.End
$spaceIns = $d.Range($spaceStart, $spaceStart)
$spaceIns.InsertAfter(" ")
$spaceRange = $d.Range($spaceStart, $spaceIns.End)
$spaceRange.Font.Size = 11

# "impedrof.github.io/curriculoOn/"
$urlStart = $spaceIns.End
$urlIns = $d.Range($urlStart, $urlStart)
$urlIns.InsertAfter("impedrof.github.io/curriculoOn/")
$urlRange = $d.Range($urlStart, $urlIns.End)
$urlRange.Font.Size = 11

# Wrap the e-mail text in a mailto: hyperlink (look it up fresh so the
# range reflects the text actually in the document right now).
$emailLinkRange = $d.Content
$emailLinkRange.Find.Execute("pedro.freitas14@fatec.sp.gov.br", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$d.Hyperlinks.Add($emailLinkRange, "mailto:pedro.freitas14@fatec.sp.gov.br") | Out-Null
$emailFix = $d.Content
$emailFix.Find.Execute("pedro.freitas14@fatec.sp.gov.br", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$emailFix.Font.Size = 11

# Wrap the new URL text in a web hyperlink.
$urlLinkRange = $d.Content
$urlLinkRange.Find.Execute("impedrof.github.io/curriculoOn/", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$d.Hyperlinks.Add($urlLinkRange, "https://impedrof.github.io/curriculoOn/") | Out-Null
$urlFix = $d.Content
$urlFix.Find.Execute("impedrof.github.io/curriculoOn/", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$urlFix.Font.Size = 11

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: it used to sit in the "Resumo Profissional"
#    paragraph (after "Estudo "); it now belongs right before "Celular:" in
#    the phone-numbers line (after the " | " separator).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$celularRng = $d.Content
$celularRng.Find.Execute("Celular:", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$celularStart = $d.Range($celularRng.Start, $celularRng.Start)
$d.Bookmarks.Add("_GoBack", $celularStart) | Out-Null
